# Release mCSD 3.9.0 with CP integrated
#
# 1) Update the Metadata sheet values (version bump, date, contact info,
#    jurisdiction, experimental flag).
# 2) Merge the two "Include" sheets (Endpoint Connect + mCSD Endpoint
#    Types) into a single sheet named "Include #0", and remove the
#    second (now redundant) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet updates
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Cells.Item(3, 2).Value  = "3.9.0"

# Writing the literal text "false" directly would make Excel store the
# cell as a real boolean (t="b"); build it via a formula instead and
# then convert that formula to its plain text result so it is kept as
# an ordinary shared string (t="s"), matching the source data.
$meta.Cells.Item(7, 2).Formula = '="fal"&"se"'
$meta.Range("B7").Copy()
$meta.Range("B7").PasteSpecial(-4163)

$meta.Cells.Item(8, 2).Value  = "2024-12-02T17:05:26-06:00"
$meta.Cells.Item(10, 2).Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$meta.Cells.Item(11, 2).Value = "null (iti@ihe.net)"
$meta.Cells.Item(12, 2).Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"
$meta.Cells.Item(13, 2).Value = "Global (Whole world)"

# ---------------------------------------------------------------------
# 2. Merge "Include from mCSD Endpoint Ty" into
#    "Include from Endpoint Connect", then rename the result "Include #0"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Include from Endpoint Connect")
$ws3 = $wb.Worksheets.Item("Include from mCSD Endpoint Ty")

# Remember the label text + the new System URI value before they get
# overwritten by the row shuffle below.
$sysUriLabel = $ws2.Cells.Item(7, 1).Value2
$newSysUri   = $ws3.Cells.Item(6, 2).Value2

# Grow ws2 from 7 to 10 rows: duplicate the formatting of the last
# (System URI) row down across the three new rows so every new cell
# keeps the same style used throughout the table.
$ws2.Range("A7:B7").Copy()
$ws2.Range("A8:B10").PasteSpecial(-4122)

# Old row 6 (blank separator row) moves down to row 9 (left blank).

# Old row 7 (System URI row) moves down to row 10, value updated to the
# URI that used to live in the (now removed) third sheet.
$ws2.Cells.Item(10, 1).Value = $sysUriLabel
$ws2.Cells.Item(10, 2).Value = $newSysUri

# Rows 6-8 get the three concepts that used to live in the third sheet.
$ws2.Cells.Item(6, 1).Value = "ihe-pdq"
$ws2.Cells.Item(6, 2).Value = ""
$ws2.Cells.Item(7, 1).Value = "ihe-pix"
$ws2.Cells.Item(7, 2).Value = ""
$ws2.Cells.Item(8, 1).Value = "ihe-mhd"
$ws2.Cells.Item(8, 2).Value = ""

# The third sheet's content has now been folded into ws2; delete it.
[void]$ws3.Delete()

# Rename the merged sheet.
$ws2.Name = "Include #0"

# Keep the same sheet selected/active as before the edit.
$meta.Activate()
